$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.175.73'
$ws.Range('E2').Value = '  +0.53%  '
# Row 3
$ws.Range('D3').Value = '1.653.57'
$ws.Range('E3').Value = '  +0.16%  '
# Row 4
$ws.Range('E4').Value = '  +0.46%  '
# Row 5
$ws.Range('D5').Value = '217.66'
$ws.Range('E5').Value = '  -0.28%  '
# Row 6
$ws.Range('E6').Value = '  +0.56%  '
# Row 7
$ws.Range('E7').Value = '  +0.44%  '
# Row 8
$ws.Range('E8').Value = '  +0.07%  '
# Row 9
$ws.Range('E9').Value = '  +0.57%  '
# Row 10
$ws.Range('D10').Value = '20.36'
$ws.Range('E10').Value = '  +0.00%  '
# Row 11
$ws.Range('D11').Value = '0.07811'
$ws.Range('E11').Value = '  +1.11%  '
# Row 12
$ws.Range('D12').Value = '4.515'
$ws.Range('E12').Value = '  +0.91%  '
# Row 13
$ws.Range('D13').Value = '1.654.95'
$ws.Range('E13').Value = '  +3.54%  '
# Row 14
$ws.Range('D14').Value = '1.881.10'
$ws.Range('E14').Value = '  +0.32%  '
# Row 15
$ws.Range('E15').Value = '  +0.63%  '
# Row 16
$ws.Range('D16').Value = '0.0₅8166'
$ws.Range('E16').Value = '  +0.82%  '
# Row 17
$ws.Range('D17').Value = '65.37'
$ws.Range('E17').Value = '  +0.48%  '
# Row 18
$ws.Range('D18').Value = '26.143.44'
$ws.Range('E18').Value = '  +0.38%  '
# Row 19
$ws.Range('E19').Value = '  +0.32%  '
# Row 20
$ws.Range('D20').Value = '4.592'
$ws.Range('E20').Value = '  +0.66%  '
# Row 21
$ws.Range('D21').Value = '190.96'
$ws.Range('E21').Value = '  -0.84%  '
# Row 22
$ws.Range('E22').Value = '  +0.34%  '
# Row 23
$ws.Range('D23').Value = '6.004'
$ws.Range('E23').Value = '  +0.21%  '
# Row 24
$ws.Range('D24').Value = '1.008'
$ws.Range('E24').Value = '  +0.52%  '
# Row 25
$ws.Range('D25').Value = '145.24'
$ws.Range('E25').Value = '  +4.15%  '
# Row 26
$ws.Range('D26').Value = '0.1220'
$ws.Range('E26').Value = '  -1.93%  '
# Row 27
$ws.Range('D27').Value = '7.193'
$ws.Range('E27').Value = '  -1.06%  '
# Row 28
$ws.Range('D28').Value = '15.98'
$ws.Range('E28').Value = '  -1.48%  '
# Row 29
$ws.Range('D29').Value = '1.474'
$ws.Range('E29').Value = '  +4.41%  '
# Row 30
$ws.Range('E30').Value = '  -3.47%  '
# Row 31
$ws.Range('D31').Value = '1.274'
$ws.Range('E31').Value = '  -0.17%  '
# Row 32
$ws.Range('D32').Value = '3.547'
$ws.Range('E32').Value = '  +1.44%  '
# Row 33
$ws.Range('D33').Value = '3.263'
$ws.Range('E33').Value = '  +0.56%  '
# Row 34
$ws.Range('D34').Value = '1.588'
$ws.Range('E34').Value = '  +3.13%  '
# Row 35
$ws.Range('E35').Value = '  +1.94%  '
# Row 36
$ws.Range('D36').Value = '2.423'
$ws.Range('E36').Value = '  +0.38%  '
# Row 37
$ws.Range('D37').Value = '0.9479'
$ws.Range('E37').Value = '  +0.35%  '
# Row 38
$ws.Range('D38').Value = '0.5746'
$ws.Range('E38').Value = '  +1.65%  '
# Row 39
$ws.Range('D39').Value = '0.01600'
$ws.Range('E39').Value = '  -0.46%  '
# Row 40
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '5.781'
$ws.Range('E40').Value = '  -1.42%  '
# Row 41
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '0.8475'
$ws.Range('E41').Value = '  -0.04%  '
# Row 42
$ws.Range('E42').Value = '  +0.49%  '
# Row 43
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.039.98'
$ws.Range('E43').Value = '  +3.31%  '
# Row 44
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = '103.82'
$ws.Range('E44').Value = '  +2.75%  '
# Row 45
$ws.Range('D45').Value = '1.793.93'
$ws.Range('E45').Value = '  +0.23%  '
# Row 46
$ws.Range('E46').Value = '  +0.11%  '
# Row 47
$ws.Range('E47').Value = '  -1.40%  '
# Row 48
$ws.Range('D48').Value = '1.006'
$ws.Range('E48').Value = '  -0.03%  '
# Row 49
$ws.Range('D49').Value = '0.4360'
$ws.Range('E49').Value = '  +1.76%  '
# Row 50
$ws.Range('D50').Value = '7.875'
$ws.Range('E50').Value = '  +0.25%  '
# Row 51
$ws.Range('E51').Value = '  +0.09%  '
